$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The line/extr naming sequence grows from 6 "lineN" + 8 "extrN" (14 rows)
# to 8 "lineN" + 8 "extrN" (16 rows): rows 8 and 9 become line7/line8, and
# every row from 10 downward shifts its extr-number down by two, with two
# brand new rows (16, 17) added at the end for extr7/extr8.

$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16

$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $false

$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# New rows 16 and 17 - match the bold/centered/bordered style used by the
# rest of column A (style index 1 in the original workbook: bold font,
# centered/top alignment, thin box border).
$ws.Range("A16").Value = 14
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("A16").Borders.LineStyle = 1
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("A17").VerticalAlignment = -4160
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
